$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = [double]"0.09872694490797113"
# Row 3
$ws.Range("B3").Value = [double]"0.001873948897131673"
$ws.Range("C3").Value = [double]"0.0006306285514183555"
$ws.Range("D3").Value = [double]"4.110995438647996"
$ws.Range("E3").Value = [double]"0.1184038137566855"
$ws.Range("F3").Value = [double]"0.0006379354962354877"
$ws.Range("G3").Value = [double]"0.003109962298027858"
$ws.Range("H3").Value = [double]"0.1006008938051028"
# Row 4
$ws.Range("B4").Value = [double]"0.005686057799308876"
$ws.Range("C4").Value = [double]"0.0009682525710397197"
$ws.Range("D4").Value = [double]"6.003186865106415"
$ws.Range("E4").Value = [double]"0.08945973012288584"
$ws.Range("F4").Value = [double]"0.003788311205471488"
$ws.Range("G4").Value = [double]"0.007583804393146265"
$ws.Range("H4").Value = [double]"0.10441300270728"
# Row 5
$ws.Range("B5").Value = [double]"0.009886247184159635"
$ws.Range("C5").Value = [double]"0.002175142466152302"
$ws.Range("D5").Value = [double]"4.24880357461066"
$ws.Range("E5").Value = [double]"0.04595304614011914"
$ws.Range("F5").Value = [double]"0.005623032150538674"
$ws.Range("G5").Value = [double]"0.01414946221778059"
$ws.Range("H5").Value = [double]"0.1086131920921308"
# Row 6
$ws.Range("B6").Value = [double]"0.007539458089537954"
$ws.Range("C6").Value = [double]"0.005183365104684434"
$ws.Range("D6").Value = [double]"4.332837221319683"
$ws.Range("E6").Value = [double]"0.1435524853106539"
$ws.Range("F6").Value = [double]"-0.00261978168260754"
$ws.Range("G6").Value = [double]"0.01769869786168345"
$ws.Range("H6").Value = [double]"0.1062664029975091"
# Row 7
$ws.Range("B7").Value = [double]"0.005773305760460595"
$ws.Range("C7").Value = [double]"0.004267657407586061"
$ws.Range("D7").Value = [double]"3.281201379418221"
$ws.Range("E7").Value = [double]"0.0652912102553153"
$ws.Range("F7").Value = [double]"-0.002591174271928169"
$ws.Range("G7").Value = [double]"0.01413778579284936"
$ws.Range("H7").Value = [double]"0.1045002506684317"
# Row 8
$ws.Range("B8").Value = [double]"0.006040934523172784"
$ws.Range("C8").Value = [double]"0.006065986850699108"
$ws.Range("D8").Value = [double]"2.497158871477062"
$ws.Range("E8").Value = [double]"0.06981263718499831"
$ws.Range("F8").Value = [double]"-0.005848224354743845"
$ws.Range("G8").Value = [double]"0.01793009340108942"
$ws.Range("H8").Value = [double]"0.1047678794311439"
# Row 9
$ws.Range("B9").Value = [double]"-0.0003021024022884141"
$ws.Range("C9").Value = [double]"0.007163037238903811"
$ws.Range("D9").Value = [double]"0.5688405218472289"
$ws.Range("E9").Value = [double]"0.1218258614801366"
$ws.Range("F9").Value = [double]"-0.01434144597758885"
$ws.Range("G9").Value = [double]"0.01373724117301202"
$ws.Range("H9").Value = [double]"0.09842484250568272"
# Row 10
$ws.Range("B10").Value = [double]"-0.09872694490797113"
$ws.Range("C10").Value = [double]"0.0004699929459654916"
$ws.Range("D10").Value = [double]"-218.8551781173092"
$ws.Range("E10").Value = [double]"0"
$ws.Range("F10").Value = [double]"-0.09964811724681735"
$ws.Range("G10").Value = [double]"-0.09780577256912494"
# Row 11
$ws.Range("B11").Value = [double]"-0.04351364602700453"
$ws.Range("C11").Value = [double]"0.0005093013262060855"
$ws.Range("D11").Value = [double]"-87.61760913489091"
$ws.Range("E11").Value = [double]"2.028958721782288e-298"
$ws.Range("F11").Value = [double]"-0.0445118616352957"
$ws.Range("G11").Value = [double]"-0.04251543041871338"
$ws.Range("H11").Value = [double]"0.05521329888096659"
# Row 12
$ws.Range("B12").Value = [double]"-0.03262303960568266"
$ws.Range("C12").Value = [double]"0.0005044976025524236"
$ws.Range("D12").Value = [double]"-66.33141707700115"
$ws.Range("E12").Value = [double]"8.827997858115449e-129"
$ws.Range("F12").Value = [double]"-0.0336118400459476"
$ws.Range("G12").Value = [double]"-0.03163423916541774"
$ws.Range("H12").Value = [double]"0.06610390530228846"
# Row 13
$ws.Range("B13").Value = [double]"-0.03019007823765087"
$ws.Range("C13").Value = [double]"0.0005049711993560494"
$ws.Range("D13").Value = [double]"-60.50887777340984"
$ws.Range("E13").Value = [double]"4.836163972969895e-112"
$ws.Range("F13").Value = [double]"-0.03117980690626443"
$ws.Range("G13").Value = [double]"-0.02920034956903732"
$ws.Range("H13").Value = [double]"0.06853686667032025"
# Row 14
$ws.Range("B14").Value = [double]"-0.02444980355673535"
$ws.Range("C14").Value = [double]"0.0004939565215179598"
$ws.Range("D14").Value = [double]"-50.83164043045711"
$ws.Range("E14").Value = [double]"9.717800205315234e-36"
$ws.Range("F14").Value = [double]"-0.02541794379668074"
$ws.Range("G14").Value = [double]"-0.02348166331678997"
$ws.Range("H14").Value = [double]"0.07427714135123577"
# Row 15
$ws.Range("B15").Value = [double]"-0.02052145900631555"
$ws.Range("C15").Value = [double]"0.0004825025512579649"
$ws.Range("D15").Value = [double]"-44.22641112593683"
$ws.Range("E15").Value = [double]"0.0003585011334396438"
$ws.Range("F15").Value = [double]"-0.0214671498074191"
$ws.Range("G15").Value = [double]"-0.01957576820521201"
$ws.Range("H15").Value = [double]"0.07820548590165557"
# Row 16
$ws.Range("B16").Value = [double]"-0.01949243626202506"
$ws.Range("C16").Value = [double]"0.0004769275978823685"
$ws.Range("D16").Value = [double]"-42.0492603604494"
$ws.Range("E16").Value = [double]"0.05127789995962657"
$ws.Range("F16").Value = [double]"-0.02042720030877247"
$ws.Range("G16").Value = [double]"-0.01855767221527765"
$ws.Range("H16").Value = [double]"0.07923450864594606"
# Row 17
$ws.Range("B17").Value = [double]"-0.0178986249898574"
$ws.Range("C17").Value = [double]"0.0004786324264362964"
$ws.Range("D17").Value = [double]"-38.432181137335"
$ws.Range("E17").Value = [double]"0.02125239229296594"
$ws.Range("F17").Value = [double]"-0.01883673044234811"
$ws.Range("G17").Value = [double]"-0.0169605195373667"
$ws.Range("H17").Value = [double]"0.08082831991811372"
# Row 18
$ws.Range("B18").Value = [double]"-0.01681913356791064"
$ws.Range("C18").Value = [double]"0.0004884467356788915"
$ws.Range("D18").Value = [double]"-36.18093884469873"
$ws.Range("E18").Value = [double]"0.002478133296051765"
$ws.Range("F18").Value = [double]"-0.01777647478097274"
$ws.Range("G18").Value = [double]"-0.01586179235484853"
$ws.Range("H18").Value = [double]"0.08190781134006048"
# Row 19
$ws.Range("B19").Value = [double]"-0.01327848011325406"
$ws.Range("C19").Value = [double]"0.0004799328953464158"
$ws.Range("D19").Value = [double]"-29.22270978643975"
$ws.Range("E19").Value = [double]"1.334544890867828e-05"
$ws.Range("F19").Value = [double]"-0.01421913446345442"
$ws.Range("G19").Value = [double]"-0.0123378257630537"
$ws.Range("H19").Value = [double]"0.08544846479471707"
# Row 20
$ws.Range("B20").Value = [double]"-0.009394060105720673"
$ws.Range("C20").Value = [double]"0.0004846535211115956"
$ws.Range("D20").Value = [double]"-20.4724675668505"
$ws.Range("E20").Value = [double]"0.0005152714835697191"
$ws.Range("F20").Value = [double]"-0.01034396672711282"
$ws.Range("G20").Value = [double]"-0.008444153484328526"
$ws.Range("H20").Value = [double]"0.08933288480225045"
# Row 21
$ws.Range("B21").Value = [double]"-0.006997350542101776"
$ws.Range("C21").Value = [double]"0.0004883682454495837"
$ws.Range("D21").Value = [double]"-15.12460776949153"
$ws.Range("E21").Value = [double]"0.0316492805442755"
$ws.Range("F21").Value = [double]"-0.007954537883738609"
$ws.Range("G21").Value = [double]"-0.006040163200464945"
$ws.Range("H21").Value = [double]"0.09172959436586935"
# Row 22
$ws.Range("B22").Value = [double]"-0.00507197237170635"
$ws.Range("C22").Value = [double]"0.0004817325938945756"
$ws.Range("D22").Value = [double]"-11.43307328946557"
$ws.Range("E22").Value = [double]"0.07562373728474472"
$ws.Range("F22").Value = [double]"-0.006016154039075219"
$ws.Range("G22").Value = [double]"-0.004127790704337481"
$ws.Range("H22").Value = [double]"0.09365497253626477"
# Row 23
$ws.Range("B23").Value = [double]"-0.005053488690693264"
$ws.Range("C23").Value = [double]"0.0004804374132996567"
$ws.Range("D23").Value = [double]"-11.15427575475242"
$ws.Range("E23").Value = [double]"0.0342329302726375"
$ws.Range("F23").Value = [double]"-0.00599513182736355"
$ws.Range("G23").Value = [double]"-0.004111845554022979"
$ws.Range("H23").Value = [double]"0.09367345621727786"
# Row 24
$ws.Range("B24").Value = [double]"-0.003988548144598819"
$ws.Range("C24").Value = [double]"0.0004768539692105781"
$ws.Range("D24").Value = [double]"-8.425731385120576"
$ws.Range("E24").Value = [double]"0.02054221608948432"
$ws.Range("F24").Value = [double]"-0.004923167851578455"
$ws.Range("G24").Value = [double]"-0.003053928437619182"
$ws.Range("H24").Value = [double]"0.09473839676337231"
# Row 25
$ws.Range("B25").Value = [double]"-0.001704809211074117"
$ws.Range("C25").Value = [double]"0.000469797061235487"
$ws.Range("D25").Value = [double]"-3.17930574399942"
$ws.Range("E25").Value = [double]"0.1863638740498571"
$ws.Range("F25").Value = [double]"-0.002625597585477827"
$ws.Range("G25").Value = [double]"-0.0007840208366704071"
$ws.Range("H25").Value = [double]"0.09702213569689701"
# Row 26
$ws.Range("B26").Value = [double]"0.008751562593373563"
$ws.Range("C26").Value = [double]"0.001697002157537079"
$ws.Range("D26").Value = [double]"5.246341887968151"
$ws.Range("E26").Value = [double]"0.08406971201847864"
$ws.Range("F26").Value = [double]"0.005425488618347853"
$ws.Range("G26").Value = [double]"0.01207763656839927"
$ws.Range("H26").Value = [double]"0.1074785075013447"
